$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 74 values (quarter ending 01-01-2021)
$ws.Range("B74").Value = 209510
$ws.Range("C74").Value = 33754
$ws.Range("E74").Value = 24767
$ws.Range("F74").Value = 2125
$ws.Range("G74").Value = 90203
$ws.Range("H74").Value = 57157
$ws.Range("I74").Value = 205567

# Add new row 75 (quarter 01-04-2021)
# Force the date-looking label to be stored as text (matches the source
# workbook's other period labels, which are all shared strings, not dates).
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "01-04-2021"
$ws.Range("A75").Style = "Normal"
$ws.Range("B75").Value = 213283
$ws.Range("C75").Value = 37860
$ws.Range("D75").Value = 1265
$ws.Range("E75").Value = 24276
$ws.Range("F75").Value = 2081
$ws.Range("G75").Value = 91116
$ws.Range("H75").Value = 56685
$ws.Range("I75").Value = 209472
